# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
#  1. Insert a new "Player Info" sheet before "ODI Batting" with the
#     player's ID / NAME / BATTING_HAND / BOWL_STYLE.
#  2. Rename the MATCH_CARD_LINK column to MATCH_CODE on both the
#     "ODI Batting" and "ODI Bowling" sheets, and replace the full
#     howstat.com scorecard URL values with just the numeric match code.

$wb = $excel.ActiveWorkbook

$urlPrefixPattern = 'http://www\.howstat\.com/cricket/Statistics/Matches/MatchScorecard_ODI\.asp\?MatchCode='

# ---------------------------------------------------------------------
# 1. Add the "Player Info" sheet in front of everything else.
#
# NOTE: worksheet handles returned/consumed here are positional, so once
# a new sheet is inserted, any previously-fetched handle now refers to
# whatever sheet currently sits at that same position. Re-fetch sheets
# by name after the insert instead of reusing old handles.
# ---------------------------------------------------------------------
$battingSheetBeforeInsert = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheetBeforeInsert)
$playerInfo.Name = "Player Info"
$playerInfo = $wb.Worksheets.Item("Player Info")

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "3680"
$playerInfo.Range("B2").Value = "Moises Constantino Henriques"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast Medium"

# ---------------------------------------------------------------------
# 2. "ODI Batting": column D, MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingLastRow = $battingSheet.Cells.Item(1, 1).Worksheet.UsedRange.Rows.Count
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $url = $cell.Value2
    if ($url) {
        $code = $url -replace $urlPrefixPattern, ''
        $cell.NumberFormat = "@"
        $cell.Value = $code
    }
}

# ---------------------------------------------------------------------
# 3. "ODI Bowling": column B, MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingLastRow = $bowlingSheet.Cells.Item(1, 1).Worksheet.UsedRange.Rows.Count
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $url = $cell.Value2
    if ($url) {
        $code = $url -replace $urlPrefixPattern, ''
        $cell.NumberFormat = "@"
        $cell.Value = $code
    }
}
